$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds/statistics values per row (matching FlashScore data refresh)

# Row 3
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.2
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = 1.36
$ws.Range("S3").Value = 6.5
$ws.Range("T3").Value = 1.11
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 2.1
$ws.Range("W3").Value = 2.38
$ws.Range("X3").Value = 1.53
$ws.Range("AA3").Value = 11
$ws.Range("AC3").Value = 26
$ws.Range("AH3").Value = 101
$ws.Range("AR3").Value = 4.57

# Row 4
$ws.Range("AR4").Value = 4.87
$ws.Range("AS4").Value = 1.16

# Row 5
$ws.Range("G5").Value = 2.7
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3.5
$ws.Range("L5").Value = 4
$ws.Range("Y5").Value = 6
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 26
$ws.Range("AK5").Value = 13
$ws.Range("AM5").Value = 34

# Row 10
$ws.Range("G10").Value = 2.7
$ws.Range("I10").Value = 2.38
$ws.Range("J10").Value = 3.25
$ws.Range("AJ10").Value = 10
$ws.Range("AN10").Value = 19

# Row 11
$ws.Range("L11").Value = 6
$ws.Range("U11").Value = 1.44
$ws.Range("V11").Value = 2.63
$ws.Range("AP11").Value = 1.64
$ws.Range("AQ11").Value = 2.22
$ws.Range("AR11").Value = 3.25
$ws.Range("AS11").Value = 1.33

# Row 14
$ws.Range("G14").Value = 2.5
$ws.Range("J14").Value = 3.4
$ws.Range("K14").Value = 1.91
$ws.Range("L14").Value = 4
$ws.Range("M14").Value = 1.11
$ws.Range("N14").Value = 6.5
$ws.Range("O14").Value = 1.53
$ws.Range("P14").Value = 2.5
$ws.Range("Q14").Value = 2.63
$ws.Range("R14").Value = 1.5
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 1.17
$ws.Range("Y14").Value = 6
$ws.Range("Z14").Value = 10
$ws.Range("AA14").Value = 11
$ws.Range("AG14").Value = 19
$ws.Range("AP14").Value = 2
$ws.Range("AQ14").Value = 1.85

# Row 15
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 15
$ws.Range("Q15").Value = 1.73
$ws.Range("R15").Value = 2.08

# Row 16
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 3.5
$ws.Range("Q16").Value = 1.98
$ws.Range("R16").Value = 1.88
$ws.Range("S16").Value = 3.4
$ws.Range("T16").Value = 1.3
$ws.Range("W16").Value = 1.75
$ws.Range("X16").Value = 2
$ws.Range("AD16").Value = 29
$ws.Range("AG16").Value = 15
$ws.Range("AH16").Value = 51
$ws.Range("AJ16").Value = 9
$ws.Range("AN16").Value = 23

# Row 18
$ws.Range("N18").Value = 8.5

# Row 19
$ws.Range("G19").Value = 2.45
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("Q19").Value = 2.08
$ws.Range("R19").Value = 1.73

# Row 20
$ws.Range("G20").Value = 1.48
$ws.Range("H20").Value = 4.5
$ws.Range("I20").Value = 5.75
$ws.Range("L20").Value = 6
$ws.Range("AE20").Value = 13
$ws.Range("AJ20").Value = 17
$ws.Range("AM20").Value = 67

# Row 21
$ws.Range("G21").Value = 2.45
$ws.Range("I21").Value = 2.75
$ws.Range("Y21").Value = 9
$ws.Range("AN21").Value = 21

# Row 22
$ws.Range("M22").Value = 1.07
$ws.Range("N22").Value = 9

# Row 23
$ws.Range("G23").Value = 1.45
$ws.Range("H23").Value = 4.2
$ws.Range("J23").Value = 2.05
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 9.5
$ws.Range("AE23").Value = 9.5
$ws.Range("AG23").Value = 23

# Row 24
$ws.Range("G24").Value = 1.85
$ws.Range("L24").Value = 5.5
$ws.Range("Q24").Value = 2.4
$ws.Range("R24").Value = 1.53
$ws.Range("AA24").Value = 9

# Row 25
$ws.Range("Q25").Value = 1.98
$ws.Range("R25").Value = 1.88
$ws.Range("S25").Value = 3.4
$ws.Range("T25").Value = 1.3

# Row 26
$ws.Range("M26").Value = 1.08
$ws.Range("N26").Value = 8
$ws.Range("O26").Value = 1.4
$ws.Range("P26").Value = 2.75
$ws.Range("Q26").Value = 2.25
$ws.Range("R26").Value = 1.62

# Row 28
$ws.Range("G28").Value = 2.27
$ws.Range("I28").Value = 2.9
$ws.Range("J28").Value = 2.77
$ws.Range("L28").Value = 3.45
$ws.Range("Y28").Value = 10.5
$ws.Range("Z28").Value = 14
$ws.Range("AB28").Value = 25
$ws.Range("AC28").Value = 16.5
$ws.Range("AD28").Value = 20
$ws.Range("AJ28").Value = 11
$ws.Range("AK28").Value = 17
$ws.Range("AL28").Value = 10.25
$ws.Range("AM28").Value = 37
$ws.Range("AN28").Value = 23
$ws.Range("AO28").Value = 26

# Row 29
$ws.Range("N29").Value = 13

# Row 32
$ws.Range("M32").Value = 1.11
$ws.Range("N32").Value = 6.5
